# #5: property aircraft done
# Fix the property_category column values that were mistakenly left as
# "land" on the 建物 (building) and 汽車 (car) sheets, and replace them
# with the correct category names.

$wb = $excel.ActiveWorkbook

# 建物 (building) sheet: column I holds property_category; rows 2-3 were
# wrongly set to "land" and should read "building".
$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2").Value = "building"
$wsBuilding.Range("I3").Value = "building"

# 汽車 (car) sheet: column H holds property_category; row 2 was wrongly
# set to "land" and should read "car".
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2").Value = "car"
